# The commit ("update scripts wuth new tpm") re-ran the NATMI ligand-receptor
# scoring pipeline with refreshed TPM input, which both changed the numeric
# scores for existing Sending/Target cluster combinations AND added the
# previously-absent "MuSCs" target-cluster rows. The sheet grows from a
# 5 (senders) x 3 (targets: FAPs, Inflammatory-Mac, Resolving-Mac) = 15-row
# block to a 5 x 4 (targets: FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac)
# = 20-row block. Row 1 (column headers) is untouched; rows 2-21 are
# rewritten in full below with the post-update values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 20,20

$data[0,0] = 'ECs'
$data[0,1] = 'App'
$data[0,2] = 'Fpr2'
$data[0,3] = 'FAPs'
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 128.4548946666667
$data[0,7] = 385.364684
$data[0,8] = 0.2815548034715028
$data[0,9] = 0.2815548034715028
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 1.109174333333333
$data[0,13] = 3.327523
$data[0,14] = 0.1199347472980627
$data[0,15] = 0.1199347472980627
$data[0,16] = 142.4788721553035
$data[0,17] = 1282.309849397732
$data[0,18] = 0.0337682042049104
$data[0,19] = 0.0337682042049104

$data[1,0] = 'ECs'
$data[1,1] = 'App'
$data[1,2] = 'Fpr2'
$data[1,3] = 'Inflammatory-Mac'
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 128.4548946666667
$data[1,7] = 385.364684
$data[1,8] = 0.2815548034715028
$data[1,9] = 0.2815548034715028
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 6.33823
$data[1,13] = 19.01469
$data[1,14] = 0.6853512477903235
$data[1,15] = 0.6853512477903234
$data[1,16] = 814.1766670231067
$data[1,17] = 7327.590003207961
$data[1,18] = 0.1929639358805538
$data[1,19] = 0.1929639358805537

$data[2,0] = 'ECs'
$data[2,1] = 'App'
$data[2,2] = 'Fpr2'
$data[2,3] = 'MuSCs'
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 128.4548946666667
$data[2,7] = 385.364684
$data[2,8] = 0.2815548034715028
$data[2,9] = 0.2815548034715028
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.09159266666666667
$data[2,13] = 0.274778
$data[2,14] = 0.009903892472889619
$data[2,15] = 0.009903892472889617
$data[2,16] = 11.76552634890578
$data[2,17] = 105.889737140152
$data[2,18] = 0.002788488498807333
$data[2,19] = 0.002788488498807332

$data[3,0] = 'ECs'
$data[3,1] = 'App'
$data[3,2] = 'Fpr2'
$data[3,3] = 'Resolving-Mac'
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 128.4548946666667
$data[3,7] = 385.364684
$data[3,8] = 0.2815548034715028
$data[3,9] = 0.2815548034715028
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 1.709151333333333
$data[3,13] = 5.127454
$data[3,14] = 0.1848101124387242
$data[3,15] = 0.1848101124387242
$data[3,16] = 219.5488544927262
$data[3,17] = 1975.939690434536
$data[3,18] = 0.05203417488723134
$data[3,19] = 0.05203417488723133

$data[4,0] = 'FAPs'
$data[4,1] = 'App'
$data[4,2] = 'Fpr2'
$data[4,3] = 'FAPs'
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 152.3944216666667
$data[4,7] = 457.183265
$data[4,8] = 0.3340268313936494
$data[4,9] = 0.3340268313936494
$data[4,10] = 2
$data[4,11] = 0.6666666666666666
$data[4,12] = 1.109174333333333
$data[4,13] = 3.327523
$data[4,14] = 0.1199347472980627
$data[4,15] = 0.1199347472980627
$data[4,16] = 169.0319810558439
$data[4,17] = 1521.287829502595
$data[4,18] = 0.04006142361396994
$data[4,19] = 0.04006142361396994

$data[5,0] = 'FAPs'
$data[5,1] = 'App'
$data[5,2] = 'Fpr2'
$data[5,3] = 'Inflammatory-Mac'
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 152.3944216666667
$data[5,7] = 457.183265
$data[5,8] = 0.3340268313936494
$data[5,9] = 0.3340268313936494
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 6.33823
$data[5,13] = 19.01469
$data[5,14] = 0.6853512477903235
$data[5,15] = 0.6853512477903234
$data[5,16] = 965.9108952403167
$data[5,17] = 8693.19805716285
$data[5,18] = 0.2289257056910856
$data[5,19] = 0.2289257056910856

$data[6,0] = 'FAPs'
$data[6,1] = 'App'
$data[6,2] = 'Fpr2'
$data[6,3] = 'MuSCs'
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 152.3944216666667
$data[6,7] = 457.183265
$data[6,8] = 0.3340268313936494
$data[6,9] = 0.3340268313936494
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.09159266666666667
$data[6,13] = 0.274778
$data[6,14] = 0.009903892472889619
$data[6,15] = 0.009903892472889617
$data[6,16] = 13.95821146557444
$data[6,17] = 125.62390319017
$data[6,18] = 0.003308165821182734
$data[6,19] = 0.003308165821182734

$data[7,0] = 'FAPs'
$data[7,1] = 'App'
$data[7,2] = 'Fpr2'
$data[7,3] = 'Resolving-Mac'
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 152.3944216666667
$data[7,7] = 457.183265
$data[7,8] = 0.3340268313936494
$data[7,9] = 0.3340268313936494
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 1.709151333333333
$data[7,13] = 5.127454
$data[7,14] = 0.1848101124387242
$data[7,15] = 0.1848101124387242
$data[7,16] = 260.4651289841456
$data[7,17] = 2344.18616085731
$data[7,18] = 0.06173153626741114
$data[7,19] = 0.06173153626741112

$data[8,0] = 'Inflammatory-Mac'
$data[8,1] = 'App'
$data[8,2] = 'Fpr2'
$data[8,3] = 'FAPs'
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 70.798157
$data[8,7] = 212.394471
$data[8,8] = 0.1551794599342134
$data[8,9] = 0.1551794599342134
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 1.109174333333333
$data[8,13] = 3.327523
$data[8,14] = 0.1199347472980627
$data[8,15] = 0.1199347472980627
$data[8,16] = 78.52749859170366
$data[8,17] = 706.7474873253329
$data[8,18] = 0.01861140931305973
$data[8,19] = 0.01861140931305973

$data[9,0] = 'Inflammatory-Mac'
$data[9,1] = 'App'
$data[9,2] = 'Fpr2'
$data[9,3] = 'Inflammatory-Mac'
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 70.798157
$data[9,7] = 212.394471
$data[9,8] = 0.1551794599342134
$data[9,9] = 0.1551794599342134
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 6.33823
$data[9,13] = 19.01469
$data[9,14] = 0.6853512477903235
$data[9,15] = 0.6853512477903234
$data[9,16] = 448.73500264211
$data[9,17] = 4038.615023778991
$data[9,18] = 0.1063524364973417
$data[9,19] = 0.1063524364973416

$data[10,0] = 'Inflammatory-Mac'
$data[10,1] = 'App'
$data[10,2] = 'Fpr2'
$data[10,3] = 'MuSCs'
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 70.798157
$data[10,7] = 212.394471
$data[10,8] = 0.1551794599342134
$data[10,9] = 0.1551794599342134
$data[10,10] = 1
$data[10,11] = 0.3333333333333333
$data[10,12] = 0.09159266666666667
$data[10,13] = 0.274778
$data[10,14] = 0.009903892472889619
$data[10,15] = 0.009903892472889617
$data[10,16] = 6.484591994715334
$data[10,17] = 58.36132795243801
$data[10,18] = 0.001536880685189532
$data[10,19] = 0.001536880685189532

$data[11,0] = 'Inflammatory-Mac'
$data[11,1] = 'App'
$data[11,2] = 'Fpr2'
$data[11,3] = 'Resolving-Mac'
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 70.798157
$data[11,7] = 212.394471
$data[11,8] = 0.1551794599342134
$data[11,9] = 0.1551794599342134
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 1.709151333333333
$data[11,13] = 5.127454
$data[11,14] = 0.1848101124387242
$data[11,15] = 0.1848101124387242
$data[11,16] = 121.0047644340927
$data[11,17] = 1089.042879906834
$data[11,18] = 0.02867873343862248
$data[11,19] = 0.02867873343862248

$data[12,0] = 'MuSCs'
$data[12,1] = 'App'
$data[12,2] = 'Fpr2'
$data[12,3] = 'FAPs'
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 20.703408
$data[12,7] = 62.110224
$data[12,8] = 0.04537891674549766
$data[12,9] = 0.04537891674549767
$data[12,10] = 2
$data[12,11] = 0.6666666666666666
$data[12,12] = 1.109174333333333
$data[12,13] = 3.327523
$data[12,14] = 0.1199347472980627
$data[12,15] = 0.1199347472980627
$data[12,16] = 22.963688766128
$data[12,17] = 206.673198895152
$data[12,18] = 0.005442508912531089
$data[12,19] = 0.005442508912531089

$data[13,0] = 'MuSCs'
$data[13,1] = 'App'
$data[13,2] = 'Fpr2'
$data[13,3] = 'Inflammatory-Mac'
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 20.703408
$data[13,7] = 62.110224
$data[13,8] = 0.04537891674549766
$data[13,9] = 0.04537891674549767
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 6.33823
$data[13,13] = 19.01469
$data[13,14] = 0.6853512477903235
$data[13,15] = 0.6853512477903234
$data[13,16] = 131.22296168784
$data[13,17] = 1181.00665519056
$data[13,18] = 0.03110049721490003
$data[13,19] = 0.03110049721490003

$data[14,0] = 'MuSCs'
$data[14,1] = 'App'
$data[14,2] = 'Fpr2'
$data[14,3] = 'MuSCs'
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 20.703408
$data[14,7] = 62.110224
$data[14,8] = 0.04537891674549766
$data[14,9] = 0.04537891674549767
$data[14,10] = 1
$data[14,11] = 0.3333333333333333
$data[14,12] = 0.09159266666666667
$data[14,13] = 0.274778
$data[14,14] = 0.009903892472889619
$data[14,15] = 0.009903892472889617
$data[14,16] = 1.896280347808
$data[14,17] = 17.066523130272
$data[14,18] = 0.000449427911983619
$data[14,19] = 0.000449427911983619

$data[15,0] = 'MuSCs'
$data[15,1] = 'App'
$data[15,2] = 'Fpr2'
$data[15,3] = 'Resolving-Mac'
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 20.703408
$data[15,7] = 62.110224
$data[15,8] = 0.04537891674549766
$data[15,9] = 0.04537891674549767
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 1.709151333333333
$data[15,13] = 5.127454
$data[15,14] = 0.1848101124387242
$data[15,15] = 0.1848101124387242
$data[15,16] = 35.38525738774401
$data[15,17] = 318.467316489696
$data[15,18] = 0.008386482706082929
$data[15,19] = 0.008386482706082929

$data[16,0] = 'Resolving-Mac'
$data[16,1] = 'App'
$data[16,2] = 'Fpr2'
$data[16,3] = 'FAPs'
$data[16,4] = 3
$data[16,5] = 1
$data[16,6] = 83.88319133333333
$data[16,7] = 251.649574
$data[16,8] = 0.1838599884551367
$data[16,9] = 0.1838599884551367
$data[16,10] = 2
$data[16,11] = 0.6666666666666666
$data[16,12] = 1.109174333333333
$data[16,13] = 3.327523
$data[16,14] = 0.1199347472980627
$data[16,15] = 0.1199347472980627
$data[16,16] = 93.04108282502243
$data[16,17] = 837.3697454252019
$data[16,18] = 0.02205120125359155
$data[16,19] = 0.02205120125359154

$data[17,0] = 'Resolving-Mac'
$data[17,1] = 'App'
$data[17,2] = 'Fpr2'
$data[17,3] = 'Inflammatory-Mac'
$data[17,4] = 3
$data[17,5] = 1
$data[17,6] = 83.88319133333333
$data[17,7] = 251.649574
$data[17,8] = 0.1838599884551367
$data[17,9] = 0.1838599884551367
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 6.33823
$data[17,13] = 19.01469
$data[17,14] = 0.6853512477903235
$data[17,15] = 0.6853512477903234
$data[17,16] = 531.6709598046733
$data[17,17] = 4785.03863824206
$data[17,18] = 0.1260086725064424
$data[17,19] = 0.1260086725064424

$data[18,0] = 'Resolving-Mac'
$data[18,1] = 'App'
$data[18,2] = 'Fpr2'
$data[18,3] = 'MuSCs'
$data[18,4] = 3
$data[18,5] = 1
$data[18,6] = 83.88319133333333
$data[18,7] = 251.649574
$data[18,8] = 0.1838599884551367
$data[18,9] = 0.1838599884551367
$data[18,10] = 1
$data[18,11] = 0.3333333333333333
$data[18,12] = 0.09159266666666667
$data[18,13] = 0.274778
$data[18,14] = 0.009903892472889619
$data[18,15] = 0.009903892472889617
$data[18,16] = 7.683085182730222
$data[18,17] = 69.147766644572
$data[18,18] = 0.0018209295557264
$data[18,19] = 0.0018209295557264

$data[19,0] = 'Resolving-Mac'
$data[19,1] = 'App'
$data[19,2] = 'Fpr2'
$data[19,3] = 'Resolving-Mac'
$data[19,4] = 3
$data[19,5] = 1
$data[19,6] = 83.88319133333333
$data[19,7] = 251.649574
$data[19,8] = 0.1838599884551367
$data[19,9] = 0.1838599884551367
$data[19,10] = 3
$data[19,11] = 1
$data[19,12] = 1.709151333333333
$data[19,13] = 5.127454
$data[19,14] = 0.1848101124387242
$data[19,15] = 0.1848101124387242
$data[19,16] = 143.3690683116218
$data[19,17] = 1290.321614804596
$data[19,18] = 0.03397918513937635
$data[19,19] = 0.03397918513937635

$ws.Range("A2:T21").Value = $data
